$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.593.12"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.15%  '

$ws.Range('D3').Value = "'2.414.38"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.47%  '

$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').Value = "'306.67"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.32%  '

$ws.Range('D6').Value = "'97.51"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.59%  '

$ws.Range('D7').Value = "'0.509"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.60%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('E9').Value = '  -1.24%  '

$ws.Range('D10').Value = "'35.22"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.27%  '

$ws.Range('D11').Value = "'0.0799"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.22%  '

$ws.Range('E12').Value = '  +2.67%  '

$ws.Range('D13').Value = "'18.52"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.52%  '

$ws.Range('D14').Value = "'6.89"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.37%  '

$ws.Range('D15').Value = "'2.780.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.19%  '

$ws.Range('D16').Value = "'2.405.53"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.23%  '

$ws.Range('D17').Value = "'0.825"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.74%  '

$ws.Range('D18').Value = "'43.579.78"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.18%  '

$ws.Range('D19').Value = "'6.43"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.56%  '

$ws.Range('D20').Value = "'12.15"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.54%  '

$ws.Range('E21').Value = '  +1.46%  '

$ws.Range('D22').Value = "'68.33"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.31%  '

$ws.Range('D23').Value = "'238.10"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.20%  '

$ws.Range('E24').Value = '  +0.83%  '

$ws.Range('D25').Value = "'2.46"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.11%  '

$ws.Range('E26').Value = '  +0.18%  '

$ws.Range('D27').Value = "'24.97"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.96%  '

$ws.Range('E28').Value = '  -0.75%  '

$ws.Range('D29').Value = "'9.42"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.30%  '

$ws.Range('D30').Value = "'32.46"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.75%  '

$ws.Range('D31').Value = "'0.117"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +15.84%  '

$ws.Range('D32').Value = "'18.45"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.31%  '

$ws.Range('D33').Value = "'5.13"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.71%  '

$ws.Range('E34').Value = '  +0.03%  '

$ws.Range('E35').Value = '  +3.53%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = "'1.89"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.25%  '

$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = "'130.32"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +26.54%  '

$ws.Range('E38').Value = '  +6.21%  '

$ws.Range('D39').Value = "'4.40"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.49%  '

$ws.Range('E40').Value = '  -1.16%  '

$ws.Range('E41').Value = '  -0.11%  '

$ws.Range('E42').Value = '  -4.56%  '

$ws.Range('D43').Value = "'1.945.88"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.14%  '

$ws.Range('E44').Value = '  +1.52%  '

$ws.Range('E45').Value = '  +1.70%  '

$ws.Range('D46').Value = "'2.84"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.70%  '

$ws.Range('D47').Value = "'9.32"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.09%  '

$ws.Range('D48').Value = "'2.635.89"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.08%  '

$ws.Range('E49').Value = '  +3.90%  '

$ws.Range('D50').Value = "'52.75"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.23%  '

$ws.Range('D51').Value = "'72.32"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.04%  '
